$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 955.6
$ws.Cells.Item(28, 9).Value = 936.6923
$ws.Cells.Item(28, 10).Value = 990.7143
$ws.Cells.Item(28, 11).Value = 936.6923
$ws.Cells.Item(28, 12).Value = 990.7143
$ws.Cells.Item(28, 13).Value = -451.6923
$ws.Cells.Item(28, 14).Value = -1960.7143
$ws.Cells.Item(40, 8).Value = 5299.909
$ws.Cells.Item(40, 10).Value = 5849.8335
$ws.Cells.Item(40, 12).Value = 5849.8335
$ws.Cells.Item(40, 14).Value = -6199.8335
$ws.Cells.Item(70, 8).Value = 9577.076999999999
$ws.Cells.Item(70, 10).Value = 14812.75
$ws.Cells.Item(70, 12).Value = 44438.25
$ws.Cells.Item(70, 14).Value = -44978.25
$ws.Cells.Item(73, 8).Value = 9577.076999999999
$ws.Cells.Item(73, 10).Value = 14812.75
$ws.Cells.Item(73, 12).Value = 44438.25
$ws.Cells.Item(73, 14).Value = -46310.25
$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 14).ClearContents()
$ws.Cells.Item(141, 8).Value = 5478.1924
$ws.Cells.Item(141, 9).Value = 4474.6816
$ws.Cells.Item(141, 11).Value = 13424.0448
$ws.Cells.Item(141, 13).Value = -8244.0448

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1702
$ws.Cells.Item(61, 9).Value = 1702
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 1702
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).Value = -1490
$ws.Cells.Item(61, 14).ClearContents()
$ws.Cells.Item(110, 8).Value = 3482.28
$ws.Cells.Item(110, 10).Value = 2749
$ws.Cells.Item(110, 12).Value = 2749
$ws.Cells.Item(110, 14).Value = -6839
$ws.Cells.Item(122, 8).Value = 3389.5
$ws.Cells.Item(122, 9).Value = 2613.7144
$ws.Cells.Item(122, 11).Value = 7841.1432
$ws.Cells.Item(122, 13).Value = -5391.1432
$ws.Cells.Item(132, 8).Value = 1690.8684
$ws.Cells.Item(132, 9).Value = 1528.4062
$ws.Cells.Item(132, 11).Value = 4585.2186
$ws.Cells.Item(132, 13).Value = -2055.2186
$ws.Cells.Item(136, 8).Value = 1702
$ws.Cells.Item(136, 9).Value = 1702
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 5106
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 13).Value = -2556
$ws.Cells.Item(136, 14).ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 2769.0435
$ws.Cells.Item(94, 9).Value = 2374.6875
$ws.Cells.Item(94, 11).Value = 2374.6875
$ws.Cells.Item(94, 13).Value = -1923.6875
$ws.Cells.Item(103, 8).Value = 105328.5
$ws.Cells.Item(103, 10).Value = 105328.5
$ws.Cells.Item(103, 12).Value = 105328.5
$ws.Cells.Item(103, 14).Value = -107672.5
$ws.Cells.Item(134, 8).Value = 2091
$ws.Cells.Item(134, 9).Value = 1577.8422
$ws.Cells.Item(134, 10).Value = 3174.3333
$ws.Cells.Item(134, 11).Value = 4733.5266
$ws.Cells.Item(134, 12).Value = 9522.999899999999
$ws.Cells.Item(134, 13).Value = -2198.5266
$ws.Cells.Item(134, 14).Value = -14592.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1786.8636
$ws.Cells.Item(58, 9).Value = 966.9286
$ws.Cells.Item(58, 10).Value = 3221.75
$ws.Cells.Item(58, 11).Value = 966.9286
$ws.Cells.Item(58, 12).Value = 3221.75
$ws.Cells.Item(58, 13).Value = -763.9286
$ws.Cells.Item(58, 14).Value = -3627.75
$ws.Cells.Item(62, 8).Value = 250005000
$ws.Cells.Item(62, 9).Value = 10000
$ws.Cells.Item(62, 11).Value = 10000
$ws.Cells.Item(62, 13).Value = -9376
$ws.Cells.Item(65, 8).Value = 250005000
$ws.Cells.Item(65, 9).Value = 10000
$ws.Cells.Item(65, 11).Value = 50000
$ws.Cells.Item(65, 13).Value = -46880
$ws.Cells.Item(123, 8).Value = 37692.31
$ws.Cells.Item(123, 10).Value = 37692.31
$ws.Cells.Item(123, 12).Value = 37692.31
$ws.Cells.Item(123, 14).Value = -47492.31
$ws.Cells.Item(136, 8).Value = 1786.8636
$ws.Cells.Item(136, 9).Value = 966.9286
$ws.Cells.Item(136, 10).Value = 3221.75
$ws.Cells.Item(136, 11).Value = 2900.7858
$ws.Cells.Item(136, 12).Value = 9665.25
$ws.Cells.Item(136, 13).Value = -350.7857999999997
$ws.Cells.Item(136, 14).Value = -14765.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(50, 8).Value = 1323.5
$ws.Cells.Item(50, 9).Value = 116.75
$ws.Cells.Item(50, 10).Value = 2128
$ws.Cells.Item(50, 11).Value = 350.25
$ws.Cells.Item(50, 12).Value = 6384
$ws.Cells.Item(50, 13).Value = 130.75
$ws.Cells.Item(50, 14).Value = -7346
$ws.Cells.Item(53, 8).Value = 1323.5
$ws.Cells.Item(53, 9).Value = 116.75
$ws.Cells.Item(53, 10).Value = 2128
$ws.Cells.Item(53, 11).Value = 350.25
$ws.Cells.Item(53, 12).Value = 6384
$ws.Cells.Item(53, 13).Value = 130.75
$ws.Cells.Item(53, 14).Value = -7346

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 13).ClearContents()
$ws.Cells.Item(44, 8).Value = 22661.166
$ws.Cells.Item(44, 9).Value = 21242.25
$ws.Cells.Item(44, 11).Value = 21242.25
$ws.Cells.Item(44, 13).Value = -20646.25
$ws.Cells.Item(47, 8).Value = 21249.5
$ws.Cells.Item(47, 10).Value = 21249.5
$ws.Cells.Item(47, 12).Value = 21249.5
$ws.Cells.Item(47, 14).Value = -22385.5
$ws.Cells.Item(48, 8).Value = 39999.668
$ws.Cells.Item(48, 9).Value = 0
$ws.Cells.Item(48, 10).Value = 39999.668
$ws.Cells.Item(48, 11).Value = 0
$ws.Cells.Item(48, 12).Value = 39999.668
$ws.Cells.Item(48, 13).ClearContents()
$ws.Cells.Item(48, 14).Value = -40969.668
$ws.Cells.Item(102, 8).Value = 3362.2856
$ws.Cells.Item(102, 9).Value = 3202.12
$ws.Cells.Item(102, 11).Value = 3202.12
$ws.Cells.Item(102, 13).Value = -1580.12
$ws.Cells.Item(122, 8).Value = 12000
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 10).Value = 12000
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 12).Value = 36000
$ws.Cells.Item(122, 13).ClearContents()
$ws.Cells.Item(122, 14).Value = -40900

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 35714980
$ws.Cells.Item(16, 9).Value = 50000570
$ws.Cells.Item(16, 10).Value = 1012.5
$ws.Cells.Item(16, 11).Value = 50000570
$ws.Cells.Item(16, 12).Value = 1012.5
$ws.Cells.Item(16, 13).Value = -50000400
$ws.Cells.Item(16, 14).Value = -1352.5
$ws.Cells.Item(68, 8).Value = 2162.3333
$ws.Cells.Item(68, 9).Value = 2162.3333
$ws.Cells.Item(68, 11).Value = 2162.3333
$ws.Cells.Item(68, 13).Value = -1413.3333
$ws.Cells.Item(71, 8).Value = 2162.3333
$ws.Cells.Item(71, 9).Value = 2162.3333
$ws.Cells.Item(71, 11).Value = 10811.6665
$ws.Cells.Item(71, 13).Value = -7067.666499999999
$ws.Cells.Item(93, 8).Value = 3072.8235
$ws.Cells.Item(93, 9).Value = 2989.5454
$ws.Cells.Item(93, 10).Value = 3225.5
$ws.Cells.Item(93, 11).Value = 2989.5454
$ws.Cells.Item(93, 12).Value = 3225.5
$ws.Cells.Item(93, 13).Value = -1741.5454
$ws.Cells.Item(93, 14).Value = -5721.5
$ws.Cells.Item(100, 8).Value = 52298.543
$ws.Cells.Item(100, 9).Value = 79144.336
$ws.Cells.Item(100, 11).Value = 79144.336
$ws.Cells.Item(100, 13).Value = -78603.336
$ws.Cells.Item(134, 8).Value = 28207.5
$ws.Cells.Item(134, 10).Value = 28207.5
$ws.Cells.Item(134, 12).Value = 28207.5
$ws.Cells.Item(134, 14).Value = -38347.5
$ws.Cells.Item(136, 8).Value = 4594.6313
$ws.Cells.Item(136, 9).Value = 4587.875
$ws.Cells.Item(136, 10).Value = 4599.5454
$ws.Cells.Item(136, 11).Value = 13763.625
$ws.Cells.Item(136, 12).Value = 13798.6362
$ws.Cells.Item(136, 13).Value = -11213.625
$ws.Cells.Item(136, 14).Value = -18898.6362

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(48, 8).Value = 40000
$ws.Cells.Item(48, 10).Value = 40000
$ws.Cells.Item(48, 12).Value = 40000
$ws.Cells.Item(48, 14).Value = -41138
$ws.Cells.Item(49, 8).Value = 30852.334
$ws.Cells.Item(49, 10).Value = 30852.334
$ws.Cells.Item(49, 12).Value = 30852.334
$ws.Cells.Item(49, 14).Value = -31312.334
$ws.Cells.Item(113, 8).Value = 568.9286
$ws.Cells.Item(113, 9).Value = 536.1818
$ws.Cells.Item(113, 11).Value = 1608.5454
$ws.Cells.Item(113, 13).Value = 561.4546
$ws.Cells.Item(117, 8).Value = 50000
$ws.Cells.Item(117, 10).Value = 50000
$ws.Cells.Item(117, 12).Value = 50000
$ws.Cells.Item(117, 14).Value = -59178
$ws.Cells.Item(125, 8).Value = 69272.73
$ws.Cells.Item(125, 9).Value = 62000
$ws.Cells.Item(125, 10).Value = 70000
$ws.Cells.Item(125, 11).Value = 62000
$ws.Cells.Item(125, 12).Value = 70000
$ws.Cells.Item(125, 13).Value = -57080
$ws.Cells.Item(125, 14).Value = -79840
